$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) column stores plain-decimal-looking values as TEXT in the
# source data (e.g. "316.55"). Excel auto-converts such text to a Number when
# assigned directly, so for those cells we force the Text number format first.

# Update Price (D) and Volume(1h) (E) columns for changed rows
$ws.Range("D2").Value = "42.350.47"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.299.86"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.55"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.68"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.00"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.966"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.29"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "2.648.21"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "2.302.61"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "42.340.48"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.13"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.10"
$ws.Range("E23").Value = "  +6.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.15"
$ws.Range("E24").Value = "  +20.39%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.77"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.84"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.30"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.136"
$ws.Range("E34").Value = "  +5.65%  "

# Rows 35-36: Kaspa and WEMIXToken swapped positions (with updated Volume(1h) values)
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.60"
$ws.Range("E35").Value = "  -10.60%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("E37").Value = "  +5.29%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.77"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.86"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.50"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "82.47"
$ws.Range("E46").Value = "  +10.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.06"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.88"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "1.589.99"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.10"
$ws.Range("E51").Value = "  -5.06%  "
